$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5", "D6", "D9", "D10", "D12", "D14", "D15", "D17", "D18", "D20", "D21", "D23", "D24", "D25", "D26", "D27", "D29", "D30", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D50", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '51.537.35'
$ws.Range('E2').Value = '  -0.20%  '
$ws.Range('D3').Value = '3.103.23'
$ws.Range('E3').Value = '  +2.15%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '385.82'
$ws.Range('E5').Value = '  +1.72%  '
$ws.Range('D6').Value = '103.77'
$ws.Range('E6').Value = '  +0.72%  '
$ws.Range('E7').Value = '  -1.17%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = '0.586'
$ws.Range('E9').Value = '  -1.37%  '
$ws.Range('D10').Value = '37.16'
$ws.Range('E10').Value = '  +0.80%  '
$ws.Range('E11').Value = '  +0.06%  '
$ws.Range('D12').Value = '0.0857'
$ws.Range('E12').Value = '  -0.47%  '
$ws.Range('D13').Value = '3.601.14'
$ws.Range('E13').Value = '  +2.26%  '
$ws.Range('D14').Value = '18.56'
$ws.Range('E14').Value = '  +0.04%  '
$ws.Range('D15').Value = '7.82'
$ws.Range('D16').Value = '3.109.55'
$ws.Range('E16').Value = '  +2.53%  '
$ws.Range('D17').Value = '0.998'
$ws.Range('E17').Value = '  +1.95%  '
$ws.Range('D18').Value = '10.87'
$ws.Range('E18').Value = '  +3.26%  '
$ws.Range('D19').Value = '51.638.41'
$ws.Range('E19').Value = '  -0.08%  '
$ws.Range('D20').Value = '3.27'
$ws.Range('E20').Value = '  +7.18%  '
$ws.Range('D21').Value = '12.48'
$ws.Range('E21').Value = '  -0.19%  '
$ws.Range('D22').Value = '0.0₃0962'
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('D23').Value = '70.05'
$ws.Range('E23').Value = '  +0.19%  '
$ws.Range('D24').Value = '266.73'
$ws.Range('E24').Value = '  -0.72%  '
$ws.Range('D25').Value = '3.18'
$ws.Range('E25').Value = '  +0.94%  '
$ws.Range('D26').Value = '8.09'
$ws.Range('E26').Value = '  -0.97%  '
$ws.Range('D27').Value = '27.03'
$ws.Range('E27').Value = '  +2.77%  '
$ws.Range('E28').Value = '  +0.02%  '
$ws.Range('D29').Value = '7.17'
$ws.Range('E29').Value = '  -5.50%  '
$ws.Range('D30').Value = '0.166'
$ws.Range('E30').Value = '  -3.79%  '
$ws.Range('E31').Value = '  -2.01%  '
$ws.Range('E32').Value = '  +1.63%  '
$ws.Range('D33').Value = '0.0480'
$ws.Range('E33').Value = '  +6.35%  '
$ws.Range('D34').Value = '35.33'
$ws.Range('E34').Value = '  +3.40%  '
$ws.Range('D35').Value = '2.07'
$ws.Range('E35').Value = '  +0.77%  '
$ws.Range('D36').Value = '50.02'
$ws.Range('E36').Value = '  -0.73%  '
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  -0.13%  '
$ws.Range('D38').Value = '3.36'
$ws.Range('E38').Value = '  +1.76%  '
$ws.Range('D39').Value = '0.292'
$ws.Range('E39').Value = '  +0.58%  '
$ws.Range('D40').Value = '1.87'
$ws.Range('E40').Value = '  +0.24%  '
$ws.Range('D41').Value = '129.15'
$ws.Range('E41').Value = '  +1.32%  '
$ws.Range('B42').Value = 'Stellar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D42').Value = '0.116'
$ws.Range('E42').Value = '  -0.20%  '
$ws.Range('B43').Value = 'Celestia'
$ws.Range('C43').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D43').Value = '16.60'
$ws.Range('E43').Value = '  -2.92%  '
$ws.Range('D44').Value = '2.51'
$ws.Range('E44').Value = '  -3.05%  '
$ws.Range('D45').Value = '3.77'
$ws.Range('E45').Value = '  +0.93%  '
$ws.Range('D46').Value = '22.16'
$ws.Range('E46').Value = '  +0.92%  '
$ws.Range('D47').Value = '2.52'
$ws.Range('E47').Value = '  +5.19%  '
$ws.Range('E48').Value = '  -2.41%  '
$ws.Range('D49').Value = '2.074.72'
$ws.Range('E49').Value = '  +1.94%  '
$ws.Range('D50').Value = '0.927'
$ws.Range('E50').Value = '  +18.05%  '
$ws.Range('D51').Value = '0.0324'
$ws.Range('E51').Value = '  +1.16%  '
